$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I (2021) mirrors the style of column H (2020) for rows 4-25.
# Copy H4:H25 formatting into I4:I25 first so the new column inherits the
# same cell styles (borders, number format, etc.) as the adjacent 2020 column.
$ws.Range("H4:H25").Copy()
$ws.Range("I4:I25").PasteSpecial(-4122)

# Now fill in the 2021 values for column I.
$ws.Range("I4").Value = 2021

$ws.Range("I5").Value = 48.5
$ws.Range("I7").Value = 48.8
$ws.Range("I8").Value = 48.2
$ws.Range("I10").Value = 58.2
$ws.Range("I11").Value = 42.4
$ws.Range("I12").Value = 40.700000000000003
$ws.Range("I14").Value = 41.5
$ws.Range("I15").Value = 52.6
$ws.Range("I17").Value = 67.099999999999994
$ws.Range("I18").Value = 62
$ws.Range("I19").Value = 46.9
$ws.Range("I20").Value = 55.8
$ws.Range("I21").Value = 42.7
$ws.Range("I22").Value = 48.3
$ws.Range("I23").Value = 39.700000000000003
$ws.Range("I24").Value = 38.1
$ws.Range("I25").Value = 44.7

# Rows 6, 9, 13 and 16 are section-header rows that stay blank in every
# year column (including the new 2021 column) - only the formatting from
# the copy above applies there, no value to set.

# Collapse the saved selection back to the default top-left cell (A1),
# matching the cleaned-up sheetView in the target workbook.
$ws.Range("A1").Select()
